# Finalized recalibration (none needed)
# Channels that were masked get "MASKED" noted, and a new drawer/channel
# entry (LBA10, chan 37, gain 1) is logged as "NOT IN UPDATE" with its date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 & 3 (LBC13) were masked channels -> note it in the Notes column (F).
$ws.Range("F2").Value = "MASKED"
$ws.Range("F3").Value = "MASKED"

# New row 4: drawer LBA10, channel 37, gain 1 - not part of this update.
$ws.Range("A4").Value = "LBA10"
$ws.Range("B4").Value = 37
$ws.Range("C4").Value = 1

# Populate the "Notes" column before the date column so the shared-string
# table ends up in the same append order as the authored workbook.
$ws.Range("F4").Value = "NOT IN UPDATE"

# "New Date" column (E) uses a text number format so "7/21" stays literal.
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "7/21"

# The extra text in row 4 wraps at the default column width, growing the
# row to fit two lines.
$ws.Rows.Item(4).RowHeight = 27

# Reflect the final cursor position left in the sheet.
$ws.Range("E5").Select()

# Page was set up for a portrait printout.
$ws.PageSetup.Orientation = 1
